# [DM] [add field] scinario_10
# Adds a new "skill_cd" field definition row (row 16) to the field-definition
# table and a new "スキル" (skill) data column (H) to the data table on the
# "装備マスタ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 16 for the new field definition -----------------
# (pushes the separator row + data table down by one row)
$ws.Rows("16:16").Insert()

# Fill the new field-definition row. Reuse the existing "整数値" cell (B9)
# formatting/shared-string via copy so the new cell points at the same
# shared string as the other "整数値" field rows.
$ws.Range("A16").Value = "skill_cd"
$ws.Range("B9").Copy()
$ws.Range("B16").PasteSpecial(-4163)
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "スキル"

# Match the field-definition row styling (border, no fill) used by sibling
# rows such as row 9/12/13.
$ws.Range("A16:F16").Borders.LineStyle = 1

# --- 2. Add the new "スキル" data column (H) ---------------------------------
$ws.Range("H19").Value = "スキル"
$ws.Range("H19").Borders.LineStyle = 1

$ws.Range("H20").Value = 2
$ws.Range("H21").Value = 2
$ws.Range("H22").Value = 2

$ws.Range("H20:H39").Borders.LineStyle = 1

# --- 3. Restore the view state ----------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("H23").Select()
